{"js": "// The document body contains a single table of 20 rows x 5 columns\n// (100 cells total), each holding one arithmetic expression such as\n// \"86-82=4\". The commit replaces every cell's expression with a new\n// one, in table (row-major) order. Duplicate \"before\" values exist\n// (e.g. \"81-8=73\" appears twice but maps to two different results),\n// so we must not do a global text search/replace -- we walk the\n// table cells in document order and overwrite each one positionally.\n\nconst newValues = [\"18+76=94\", \"8+9=17\", \"46-5=41\", \"63+36=99\", \"47-11=36\", \"14+20=34\", \"92-18=74\", \"54-32=22\", \"36+57=93\", \"22+49=71\", \"96-47=49\", \"50+48=98\", \"30+53=83\", \"80+0=80\", \"6+53=59\", \"7+15=22\", \"74-56=18\", \"80-9=71\", \"8+34=42\", \"87-9=78\", \"60+18=78\", \"29+65=94\", \"82-60=22\", \"35-27=8\", \"80-47=33\", \"68+9=77\", \"87-36=51\", \"87-40=47\", \"75-27=48\", \"10+53=63\", \"79-74=5\", \"22+57=79\", \"29+58=87\", \"67-15=52\", \"75-9=66\", \"28-6=22\", \"74-27=47\", \"67-64=3\", \"70-16=54\", \"67-16=51\", \"23+19=42\", \"12+72=84\", \"81-43=38\", \"85-7=78\", \"49+50=99\", \"64-29=35\", \"18-12=6\", \"11+42=53\", \"81-18=63\", \"64+1=65\", \"76+20=96\", \"72-3=69\", \"68-32=36\", \"52+37=89\", \"91-11=80\", \"90-13=77\", \"65-16=49\", \"31-11=20\", \"97-42=55\", \"9+86=95\", \"6+5=11\", \"58-0=58\", \"21+63=84\", \"34+47=81\", \"44-27=17\", \"44-39=5\", \"78+17=95\", \"78-47=31\", \"48+26=74\", \"98-2=96\", \"75+5=80\", \"71+6=77\", \"55-32=23\", \"40-37=3\", \"7+74=81\", \"35-0=35\", \"44+34=78\", \"60-58=2\", \"72+14=86\", \"32-2=30\", \"80-51=29\", \"46-26=20\", \"12+47=59\", \"42+57=99\", \"45+15=60\", \"12+7=19\", \"45+15=60\", \"32+60=92\", \"77+7=84\", \"15+44=59\", \"60-32=28\", \"29+63=92\", \"54-13=41\", \"13+18=31\", \"44+55=99\", \"96-34=62\", \"85-35=50\", \"40+39=79\", \"51+19=70\", \"98-84=14\"];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Load each row's cells collection.\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nlet i = 0;\nfor (const row of rows.items) {\n  for (const cell of row.cells.items) {\n    if (i >= newValues.length) break;\n    cell.value = newValues[i];\n    i++;\n  }\n}\nawait context.sync();\n", "ps1": "# The document body contains a single table of 20 rows x 5 columns\n# (100 cells total), each holding one arithmetic expression such as\n# \"86-82=4\". The commit replaces every cell's expression with a new\n# one, in table (row-major) order. Duplicate \"before\" values exist\n# (e.g. \"81-8=73\" appears twice but maps to two different results),\n# so we must not do a global Find/Replace across the document -- we\n# walk the table cells in order (Cell(row, col), row-major) and\n# overwrite each one positionally.\n\n$newValues = @(\n  \"18+76=94\", \"8+9=17\", \"46-5=41\", \"63+36=99\", \"47-11=36\",\n  \"14+20=34\", \"92-18=74\", \"54-32=22\", \"36+57=93\", \"22+49=71\",\n  \"96-47=49\", \"50+48=98\", \"30+53=83\", \"80+0=80\", \"6+53=59\",\n  \"7+15=22\", \"74-56=18\", \"80-9=71\", \"8+34=42\", \"87-9=78\",\n  \"60+18=78\", \"29+65=94\", \"82-60=22\", \"35-27=8\", \"80-47=33\",\n  \"68+9=77\", \"87-36=51\", \"87-40=47\", \"75-27=48\", \"10+53=63\",\n  \"79-74=5\", \"22+57=79\", \"29+58=87\", \"67-15=52\", \"75-9=66\",\n  \"28-6=22\", \"74-27=47\", \"67-64=3\", \"70-16=54\", \"67-16=51\",\n  \"23+19=42\", \"12+72=84\", \"81-43=38\", \"85-7=78\", \"49+50=99\",\n  \"64-29=35\", \"18-12=6\", \"11+42=53\", \"81-18=63\", \"64+1=65\",\n  \"76+20=96\", \"72-3=69\", \"68-32=36\", \"52+37=89\", \"91-11=80\",\n  \"90-13=77\", \"65-16=49\", \"31-11=20\", \"97-42=55\", \"9+86=95\",\n  \"6+5=11\", \"58-0=58\", \"21+63=84\", \"34+47=81\", \"44-27=17\",\n  \"44-39=5\", \"78+17=95\", \"78-47=31\", \"48+26=74\", \"98-2=96\",\n  \"75+5=80\", \"71+6=77\", \"55-32=23\", \"40-37=3\", \"7+74=81\",\n  \"35-0=35\", \"44+34=78\", \"60-58=2\", \"72+14=86\", \"32-2=30\",\n  \"80-51=29\", \"46-26=20\", \"12+47=59\", \"42+57=99\", \"45+15=60\",\n  \"12+7=19\", \"45+15=60\", \"32+60=92\", \"77+7=84\", \"15+44=59\",\n  \"60-32=28\", \"29+63=92\", \"54-13=41\", \"13+18=31\", \"44+55=99\",\n  \"96-34=62\", \"85-35=50\", \"40+39=79\", \"51+19=70\", \"98-84=14\"\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\n\n$i = 0\nfor ($r = 1; $r -le $rows; $r++) {\n  for ($c = 1; $c -le $cols; $c++) {\n    $cell = $t.Cell($r, $c)\n    $cell.Range.Text = $newValues[$i]\n    $i++\n  }\n}\n"}
